# Etapa 2: elimina los datos antiguos (paises "Bolivia", "China", ... en texto
# completo y encabezados "anyo"/"pais"/"tasa_homicidios") y los reemplaza con
# el nuevo layout (Year / Nationality code / Homicide Rate, codigos ISO3) mas
# las filas agregadas para Cuba, Peru, Rusia y los anios adicionales de China,
# Ucrania y Senegal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: seed the shared-string table in the exact order the new
# strings must first appear (the final sheet layout below does NOT visit the
# values in this order, since the rows end up grouped/re-sorted by country).
# Writing them first, off to the side, pins their shared-string indices; the
# scratch cells are cleared again at the very end.
$stringOrder = @(
    "Year", "Nationality code", "Homicide Rate",
    "BOL", "CHN", "SEN", "PRT", "VEN", "UKR", "CUB", "PER", "RUS"
)
for ($i = 0; $i -lt $stringOrder.Length; $i++) {
    $ws.Cells.Item(200, $i + 1).Value = $stringOrder[$i]
}

# --- Phase 2: write the final A1:C29 layout.
$ws.Cells.Item(1,1).Value = "Year"
$ws.Cells.Item(1,2).Value = "Nationality code"
$ws.Cells.Item(1,3).Value = "Homicide Rate"

$ws.Cells.Item(2,1).Value = 2008
$ws.Cells.Item(2,2).Value = "BOL"
$ws.Cells.Item(2,3).Value = 8.52

$ws.Cells.Item(3,1).Value = 2009
$ws.Cells.Item(3,2).Value = "BOL"
$ws.Cells.Item(3,3).Value = 8.2899999999999991

$ws.Cells.Item(4,1).Value = 2010
$ws.Cells.Item(4,2).Value = "BOL"
$ws.Cells.Item(4,3).Value = 12.67

$ws.Cells.Item(5,1).Value = 2011
$ws.Cells.Item(5,2).Value = "BOL"
$ws.Cells.Item(5,3).Value = 12.1

$ws.Cells.Item(6,1).Value = 2012
$ws.Cells.Item(6,2).Value = "BOL"
$ws.Cells.Item(6,3).Value = 11.77

$ws.Cells.Item(7,1).Value = 2013
$ws.Cells.Item(7,2).Value = "BOL"
$ws.Cells.Item(7,3).Value = 8.84

$ws.Cells.Item(8,1).Value = 2014
$ws.Cells.Item(8,2).Value = "BOL"
$ws.Cells.Item(8,3).Value = 8.8000000000000007

$ws.Cells.Item(9,1).Value = 2021
$ws.Cells.Item(9,2).Value = "CHN"
$ws.Cells.Item(9,3).Value = 0.45

$ws.Cells.Item(10,1).Value = 2022
$ws.Cells.Item(10,2).Value = "CHN"
$ws.Cells.Item(10,3).Value = 0.4

$ws.Cells.Item(11,1).Value = 2008
$ws.Cells.Item(11,2).Value = "PRT"
$ws.Cells.Item(11,3).Value = 1.17

$ws.Cells.Item(12,1).Value = 2009
$ws.Cells.Item(12,2).Value = "PRT"
$ws.Cells.Item(12,3).Value = 1.2

$ws.Cells.Item(13,1).Value = 2022
$ws.Cells.Item(13,2).Value = "PRT"
$ws.Cells.Item(13,3).Value = 0.72

$ws.Cells.Item(14,1).Value = 2013
$ws.Cells.Item(14,2).Value = "VEN"
$ws.Cells.Item(14,3).Value = 79

$ws.Cells.Item(15,1).Value = 2018
$ws.Cells.Item(15,2).Value = "VEN"
$ws.Cells.Item(15,3).Value = 36.69

$ws.Cells.Item(16,1).Value = 2020
$ws.Cells.Item(16,2).Value = "CUB"
$ws.Cells.Item(16,3).Value = 4.38

$ws.Cells.Item(17,1).Value = 2021
$ws.Cells.Item(17,2).Value = "CUB"
$ws.Cells.Item(17,3).Value = 4.3

$ws.Cells.Item(18,1).Value = 2022
$ws.Cells.Item(18,2).Value = "CUB"
$ws.Cells.Item(18,3).Value = 4.34

$ws.Cells.Item(19,1).Value = 2008
$ws.Cells.Item(19,2).Value = "PER"
$ws.Cells.Item(19,3).Value = 5.27

$ws.Cells.Item(20,1).Value = 2009
$ws.Cells.Item(20,2).Value = "PER"
$ws.Cells.Item(20,3).Value = 5.38

$ws.Cells.Item(21,1).Value = 2010
$ws.Cells.Item(21,2).Value = "PER"
$ws.Cells.Item(21,3).Value = 5.4

$ws.Cells.Item(22,1).Value = 2022
$ws.Cells.Item(22,2).Value = "PER"
$ws.Cells.Item(22,3).Value = 7.08

$ws.Cells.Item(23,1).Value = 2022
$ws.Cells.Item(23,2).Value = "RUS"
$ws.Cells.Item(23,3).Value = 7.07

$ws.Cells.Item(24,1).Value = 2011
$ws.Cells.Item(24,2).Value = "UKR"
$ws.Cells.Item(24,3).Value = 4.88

$ws.Cells.Item(25,1).Value = 2013
$ws.Cells.Item(25,2).Value = "UKR"
$ws.Cells.Item(25,3).Value = 5.72

$ws.Cells.Item(26,1).Value = 2015
$ws.Cells.Item(26,2).Value = "UKR"
$ws.Cells.Item(26,3).Value = 6.16

$ws.Cells.Item(27,1).Value = 2016
$ws.Cells.Item(27,2).Value = "UKR"
$ws.Cells.Item(27,3).Value = 6.16

$ws.Cells.Item(28,1).Value = 2022
$ws.Cells.Item(28,2).Value = "UKR"
$ws.Cells.Item(28,3).Value = 3.86

$ws.Cells.Item(29,1).Value = 2015
$ws.Cells.Item(29,2).Value = "SEN"
$ws.Cells.Item(29,3).Value = 0.27

# --- Phase 3: clear the scratch cells used to pin shared-string order.
$ws.Range("A200:L200").Clear()

# --- Column layout: column B ("Nationality code") gets an explicit width,
# matching the new second column of data.
$ws.Columns.Item(2).ColumnWidth = 21.86

# --- A9:C10 (the new China rows) pick up an explicit (no-op) fill format.
$ws.Range("A9:C10").Interior.ColorIndex = -4142

# --- Selection moves to F27, matching the author's last cursor position.
$ws.Range("F27").Select()
